# Regenerate merged AHB files
#
# 1. Rename the header labels in row 1 (A1:U1):
#       "<Name>_old" -> "<Name>_FV2410"
#       "<Name>_new" -> "<Name>_FV2504"
# 2. Turn the used range into an Excel Table ("Table1") with an AutoFilter.
# 3. Freeze the header row (split after row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count()
$lastCol = $used.Columns.Count()

# --- 1. Rename the header row values -------------------------------------
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = [string]$cell.Value()
    if ($val.EndsWith("_old")) {
        $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2410"
    } elseif ($val.EndsWith("_new")) {
        $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2504"
    }
}

# --- 2. Create the table over the full used range -------------------------
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# --- 3. Freeze the header row ---------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
